$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = -20.48189999999998
$ws.Range("A12").Value = -22.51500000000003
$ws.Range("E13").Value = 12.7966
$ws.Range("A18").Value = -22.40690000000002
